# The timesheet's "Hours" column (C) was missing an entry for the Wed
# of the first week (row 3). Fill it in the same way the other days in
# that column are recorded (4 hours), which is what the author forgot
# to do before committing ("Forgot to update timesheet").
#
# The weekly/overall SUM() formulas in C8 and C46 recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = 4

# Leave the cursor where the author would have landed after typing the
# value into C3 and pressing Enter.
$null = $ws.Range("C4").Select()
